# Apply the commit's change:
#  1. In the "ODI Batting" sheet, the (already-empty) INNING_NUMBER cells
#     B8, B9 and B11 are removed entirely (cleared).
#  2. A new worksheet "ODI Batting Extra" (sheetId 4) is appended after
#     "ODI Bowling", containing MATCH_CODE / BATTING_POSITION / NUM_4 /
#     NUM_6 / PERCENT_RUNS_OF_TOTAL / MAN_OF_MATCH columns.

$wb = $excel.ActiveWorkbook

# --- 1. Clear the stray empty cells in "ODI Batting" ----------------------
$odiBatting = $wb.Worksheets.Item("ODI Batting")
$odiBatting.Range("B8").Clear()
$odiBatting.Range("B9").Clear()
$odiBatting.Range("B11").Clear()

# --- 2. Add the new "ODI Batting Extra" worksheet at the end --------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$extra = $wb.Worksheets.Add($null, $lastSheet)
$extra.Name = "ODI Batting Extra"

# Header row text.
$headers = @("MATCH_CODE", "BATTING_POSITION", "NUM_4", "NUM_6", "PERCENT_RUNS_OF_TOTAL", "MAN_OF_MATCH")
for ($col = 1; $col -le $headers.Length; $col++) {
    $extra.Cells.Item(1, $col).Value = $headers[$col - 1]
}
# Match the bold / centered / top-aligned / thin-bordered header style
# already used on the other sheets by copying its format.
$odiBatting.Range("A1:F1").Copy()
$extra.Range("A1:F1").PasteSpecial(-4122)  # xlPasteFormats

# MATCH_CODE (column A) reuses the exact same match-code text values that
# already live in "ODI Batting"!D2:D12 (same order), so copy them across
# instead of re-typing - this keeps them as real text instead of numbers.
$odiBatting.Range("D2:D12").Copy()
$extra.Range("A2:A12").PasteSpecial(-4104)  # xlPasteAll
$excel.CutCopyMode = $false

# BATTING_POSITION: numeric when present, otherwise blank.
$battingPositions = @($null, 9, 8, $null, 8, 8, 8, $null, 9, 8, 8)
for ($i = 0; $i -lt $battingPositions.Length; $i++) {
    $cell = $extra.Cells.Item($i + 2, 2)
    if ($null -eq $battingPositions[$i]) {
        $cell.Value = "'"
    } else {
        $cell.Value = $battingPositions[$i]
    }
}

# NUM_4, NUM_6, PERCENT_RUNS_OF_TOTAL: text columns, possibly blank. Force
# text (instead of letting Excel coerce numeric-/percent-looking strings
# into numbers) with a leading apostrophe, the same trick Excel's UI uses.
$num4 = @("", "0", "2", "", "6", "1", "", "0", "0", "", "2")
$num6 = @("", "0", "0", "", "1", "0", "", "0", "0", "", "1")
$percentOfTotal = @("", "1.38%", "4.47%", "", "23.33%", "3.59%", "", "", "", "", "9.67%")

for ($i = 0; $i -lt $num4.Length; $i++) {
    $r = $i + 2
    $extra.Cells.Item($r, 3).Value = "'" + $num4[$i]
    $extra.Cells.Item($r, 4).Value = "'" + $num6[$i]
    $extra.Cells.Item($r, 5).Value = "'" + $percentOfTotal[$i]
}

# MAN_OF_MATCH: plain text, always "NO" for these rows.
for ($r = 2; $r -le 12; $r++) {
    $extra.Cells.Item($r, 6).Value = "NO"
}

$extra.Range("A1").Select()
